# Applies the MeshSr_OFTest.xlsx commit:
#   "add testCase and new python test files"
#
# - Adds TEST SUITE 80 block (rows 108-138) with TestCase 80.10 .. 80.300
# - Two cells ("1.3 not support") and four cells (rich-text note with a
#   red "???????" run) are added in column D for some of those rows
# - Widens column D to fit the new long text
# - Moves the active selection to the new last entered cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Register the red "Arial Narrow" 11pt font used by the rich-text run
#    further down, by touching a scratch cell and clearing it again so
#    the font lands in the shared font table without leaving data behind.
# ---------------------------------------------------------------------
$ws.Range("ZZ1").Value = "x"
$ws.Range("ZZ1").Font.Size = 11
$ws.Range("ZZ1").Font.Color = 255
$ws.Range("ZZ1").Font.Name = "Arial Narrow"
$ws.Range("ZZ1").Clear()

# ---------------------------------------------------------------------
# 2. Widen column D for the new, much longer comments.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 59.785714285714285

# ---------------------------------------------------------------------
# 3. TEST SUITE 80 header (bold, 12pt, Arial Narrow - same look as the
#    other "TEST SUITE" rows already in the sheet).
# ---------------------------------------------------------------------
$ws.Range("A108").Value = "TEST SUITE 80"
$ws.Range("A108").Font.Bold = $true
$ws.Range("A108").Font.Size = 12
$ws.Range("A108").Font.Name = "Arial Narrow"
$ws.Range("A108").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 4. Plain TestCase rows (column A only) - these simply inherit the
#    column's default style, matching the existing TestCase rows.
# ---------------------------------------------------------------------
$ws.Range("A109").Value = "TestCase 80.10"
$ws.Range("A110").Value = "TestCase 80.20"
$ws.Range("A111").Value = "TestCase 80.30"
$ws.Range("A112").Value = "TestCase 80.40"
$ws.Range("A113").Value = "TestCase 80.50"
$ws.Range("A114").Value = "TestCase 80.60"
$ws.Range("A115").Value = "TestCase 80.70"
$ws.Range("A116").Value = "TestCase 80.80"
$ws.Range("A117").Value = "TestCase 80.90"
$ws.Range("A118").Value = "TestCase 80.100"
$ws.Range("A119").Value = "TestCase 80.110"
$ws.Range("A120").Value = "TestCase 80.120"
$ws.Range("A121").Value = "TestCase 80.130"

# ---------------------------------------------------------------------
# 5. TestCase rows with a "1.3 not support" note in column D.
# ---------------------------------------------------------------------
$ws.Range("A122").Value = "TestCase 80.140"
$ws.Range("D122").Value = "1.3 not support"

$ws.Range("A123").Value = "TestCase 80.150"
$ws.Range("D123").Value = "1.3 not support"

$ws.Range("A124").Value = "TestCase 80.160"

$ws.Range("A125").Value = "TestCase 80.170"
$ws.Range("D125").Value = "1.3 not support"

$ws.Range("A126").Value = "TestCase 80.180"
$ws.Range("D126").Value = "1.3 not support"

$ws.Range("A127").Value = "TestCase 80.190"
$ws.Range("D127").Value = "1.3 not support"

$ws.Range("A128").Value = "TestCase 80.200"
$ws.Range("D128").Value = "1.3 not support"

# ---------------------------------------------------------------------
# 6. TestCase rows with the rich-text comment: a plain sentence followed
#    by a red "???????" run in Arial Narrow 11pt.
# ---------------------------------------------------------------------
$richText = "verify the handling of fragments is consistent with the returned configuration " + "???????"

$ws.Range("A129").Value = "TestCase 80.210"
$ws.Range("D129").Value = $richText
$ws.Range("D129").Characters(80, 7).Font.Color = 255
$ws.Range("D129").Characters(80, 7).Font.Name = "Arial Narrow"
$ws.Range("D129").Characters(80, 7).Font.Size = 11

$ws.Range("A130").Value = "TestCase 80.220"
$ws.Range("D130").Value = $richText
$ws.Range("D130").Characters(80, 7).Font.Color = 255
$ws.Range("D130").Characters(80, 7).Font.Name = "Arial Narrow"
$ws.Range("D130").Characters(80, 7).Font.Size = 11

$ws.Range("A131").Value = "TestCase 80.230"
$ws.Range("D131").Value = $richText
$ws.Range("D131").Characters(80, 7).Font.Color = 255
$ws.Range("D131").Characters(80, 7).Font.Name = "Arial Narrow"
$ws.Range("D131").Characters(80, 7).Font.Size = 11

$ws.Range("A132").Value = "TestCase 80.240"
$ws.Range("D132").Value = $richText
$ws.Range("D132").Characters(80, 7).Font.Color = 255
$ws.Range("D132").Characters(80, 7).Font.Name = "Arial Narrow"
$ws.Range("D132").Characters(80, 7).Font.Size = 11

# ---------------------------------------------------------------------
# 7. Remaining plain TestCase rows.
# ---------------------------------------------------------------------
$ws.Range("A133").Value = "TestCase 80.250"
$ws.Range("A134").Value = "TestCase 80.260"
$ws.Range("A135").Value = "TestCase 80.270"
$ws.Range("A136").Value = "TestCase 80.280"
$ws.Range("A137").Value = "TestCase 80.290"
$ws.Range("A138").Value = "TestCase 80.300"

# ---------------------------------------------------------------------
# 8. Leave the selection where the author left it.
# ---------------------------------------------------------------------
$ws.Range("B134").Select()

Write-Host "MeshSr_OFTest TEST SUITE 80 block added"
